$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.003055490664580075
$ws.Range("G2").Value = 17
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.002846270997462597
$ws.Range("G3").Value = 15
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.004599586612760704
$ws.Range("G4").Value = 15
$ws.Range("C5").Value = 87
$ws.Range("E5").Value = 0.004380600676719868
$ws.Range("G5").Value = 14
$ws.Range("C6").Value = 89
$ws.Range("E6").Value = 0.0027174320259105
$ws.Range("G6").Value = 14

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.003055490664580075
$ws.Range("G2").Value = 45
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.002846270997462597
$ws.Range("G3").Value = 40
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.004599586612760704
$ws.Range("G4").Value = 37
$ws.Range("C5").Value = 87
$ws.Range("E5").Value = 0.004380600676719868
$ws.Range("G5").Value = 33
$ws.Range("C6").Value = 89
$ws.Range("E6").Value = 0.0027174320259105
$ws.Range("G6").Value = 44

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.003055490664580075
$ws.Range("G2").Value = 63
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.002846270997462597
$ws.Range("G3").Value = 62
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.004599586612760704
$ws.Range("G4").Value = 59
$ws.Range("C5").Value = 87
$ws.Range("E5").Value = 0.004380600676719868
$ws.Range("G5").Value = 61
$ws.Range("C6").Value = 89
$ws.Range("E6").Value = 0.0027174320259105
$ws.Range("G6").Value = 63

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("C2").Value = 88
$ws.Range("E2").Value = 0.003055490664580075
$ws.Range("G2").Value = 74
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.002846270997462597
$ws.Range("G3").Value = 74
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.004599586612760704
$ws.Range("G4").Value = 72
$ws.Range("C5").Value = 87
$ws.Range("E5").Value = 0.004380600676719868
$ws.Range("G5").Value = 74
$ws.Range("C6").Value = 89
$ws.Range("E6").Value = 0.0027174320259105
$ws.Range("G6").Value = 74
